# Regenerate the "K" column (column G) of the save-data sheet.
# The original analysis derived K (formerly "Strike#") per row; this
# replays the recalculated values that resulted from that regen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 2
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 2
    16 = 1
    17 = 1
    18 = 2
    19 = 1
    20 = 0
    21 = 1
    22 = 1
    23 = 3
    24 = 2
    25 = 2
    26 = 2
    27 = 1
    28 = 1
    29 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
